$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells whose new values look like pure numbers remain stored as text,
# matching the original inline-string (text) representation in the workbook.
$textFormatCells = @("D5", "D6", "D10", "D16", "D20", "D21", "D22", "D25", "D26", "D30", "D31", "D33", "D34", "D35", "D37", "D38", "D39", "D41", "D42", "D43", "D46", "D47", "D49", "D50")
foreach ($cellAddr in $textFormatCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Update Price (column D) values
$ws.Range("D2").Value = "59.102.82"
$ws.Range("D3").Value = "2.571.31"
$ws.Range("D5").Value = "554.62"
$ws.Range("D6").Value = "141.23"
$ws.Range("D9").Value = "2.575.21"
$ws.Range("D10").Value = "6.64"
$ws.Range("D14").Value = "3.025.14"
$ws.Range("D15").Value = "59.077.79"
$ws.Range("D16").Value = "22.92"
$ws.Range("D18").Value = "2.574.63"
$ws.Range("D20").Value = "335.97"
$ws.Range("D21").Value = "10.29"
$ws.Range("D22").Value = "6.38"
$ws.Range("D25").Value = "0.469"
$ws.Range("D26").Value = "1.00"
$ws.Range("D29").Value = "0.0₃0768"
$ws.Range("D30").Value = "0.998"
$ws.Range("D31").Value = "6.15"
$ws.Range("D33").Value = "157.64"
$ws.Range("D34").Value = "18.96"
$ws.Range("D35").Value = "4.01"
$ws.Range("D37").Value = "0.892"
$ws.Range("D38").Value = "37.24"
$ws.Range("D39").Value = "0.852"
$ws.Range("D41").Value = "3.65"
$ws.Range("D42").Value = "291.00"
$ws.Range("D43").Value = "135.12"
$ws.Range("D46").Value = "0.590"
$ws.Range("D47").Value = "10.67"
$ws.Range("D49").Value = "0.0232"
$ws.Range("D50").Value = "18.55"
$ws.Range("D51").Value = "1.939.36"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("E5").Value = "  -2.41%  "
$ws.Range("E6").Value = "  -3.54%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  +12.28%  "
$ws.Range("E14").Value = "  -2.77%  "
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("E16").Value = "  +3.81%  "
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("E18").Value = "  -2.70%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  -2.33%  "
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("E24").Value = "  -5.33%  "
$ws.Range("E25").Value = "  +6.53%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  -3.79%  "
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("E29").Value = "  -5.66%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E32").Value = "  -3.38%  "
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("E35").Value = "  -2.49%  "
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("E39").Value = "  -5.57%  "
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("E43").Value = "  +4.66%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E48").Value = "  -3.21%  "
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("E51").Value = "  -1.15%  "
